$wb = $excel.ActiveWorkbook

# "Set Values Here" sheet: update the carbon-tax-revenue weighting row (row 8)
# so that C8 and E8 go from 0 to 5 (spreading the weight across Household Taxes
# and Payroll Taxes in addition to the existing Deficit Spending weight).
$svh = $wb.Worksheets.Item("Set Values Here")
$svh.Range("C8").Value = 5
$svh.Range("E8").Value = 5
$svh.Range("C9").Select()

# "GRA-carbontax" sheet: selection moved to B5 in the saved file.
$carbontax = $wb.Worksheets.Item("GRA-carbontax")
$carbontax.Range("B5").Select()

# Restore "About" as the active/displayed sheet, as in the original workbook.
$about = $wb.Worksheets.Item("About")
$about.Activate()

$wb.Application.Calculate()
